$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 410/411; existing rows 410-442 shift down to 412-444.
$ws.Range("A410:A411").EntireRow.Insert()

# New row 410 data
$ws.Cells.Item(410, 1).Value = 10
$ws.Cells.Item(410, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(410, 3).Value = "La Araucanía"
$ws.Cells.Item(410, 4).Value = 44826
$ws.Cells.Item(410, 5).Value = 9
$ws.Cells.Item(410, 6).Value = 100112040
$ws.Cells.Item(410, 7).Value = "Cilantro"
$ws.Cells.Item(410, 8).Value = "Sin especificar"
$ws.Cells.Item(410, 9).Value = "Primera"
$ws.Cells.Item(410, 10).Value = 110
$ws.Cells.Item(410, 11).Value = 5000
$ws.Cells.Item(410, 12).Value = 6000
$ws.Cells.Item(410, 13).Value = 5455
$ws.Cells.Item(410, 14).Value = "$/docena de atados (2 kilos)"
$ws.Cells.Item(410, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(410, 16).Value = 2728
$ws.Cells.Item(410, 17).Value = 2
$ws.Cells.Item(410, 18).Value = "Hortaliza"

# New row 411 data
$ws.Cells.Item(411, 1).Value = 10
$ws.Cells.Item(411, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(411, 3).Value = "La Araucanía"
$ws.Cells.Item(411, 4).Value = 44826
$ws.Cells.Item(411, 5).Value = 9
$ws.Cells.Item(411, 6).Value = 100112040
$ws.Cells.Item(411, 7).Value = "Cilantro"
$ws.Cells.Item(411, 8).Value = "Sin especificar"
$ws.Cells.Item(411, 9).Value = "Primera"
$ws.Cells.Item(411, 10).Value = 150
$ws.Cells.Item(411, 11).Value = 3300
$ws.Cells.Item(411, 12).Value = 3300
$ws.Cells.Item(411, 13).Value = 3300
$ws.Cells.Item(411, 14).Value = "$/docena de atados (2 kilos)"
$ws.Cells.Item(411, 15).Value = "Región Metropolitana"
$ws.Cells.Item(411, 16).Value = 1650
$ws.Cells.Item(411, 17).Value = 2
$ws.Cells.Item(411, 18).Value = "Hortaliza"
